# Summer 2024 Working Hours.xlsx — add the 2024-05-14 entry.
#
# "trying out random forests and GBM on basic math/stat courses"
#
# Appends a new log row (date, hours worked, notes) to the bottom of the
# table on Sheet1, mirroring the format of the existing rows, and leaves
# the selection on the new notes cell the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: 2024-05-14, 2.5 hours, notes about the work done.
$ws.Range("A13").Value = 45426
# Match the date formatting used by the rest of the Date column (d-mmm).
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat
$ws.Range("B13").Value = 2.5
$ws.Range("D13").Value = "trying out random forests, missForest, and GBM"

# The wrapped note in D12 now renders a touch taller under the refreshed
# default row height/font metrics.
$ws.Rows.Item(12).RowHeight = 45

# Leave the selection on the newly typed note, same as the author did.
$ws.Range("D13").Select() | Out-Null
